$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.569.49'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '1.884.42'
$ws.Range('E3').Value = '  +1.50%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.14%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '247.46'
$c.ClearFormats()
$ws.Range('E5').Value = '  +5.97%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9993'
$c.ClearFormats()
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4758'
$c.ClearFormats()
$ws.Range('E7').Value = '  +1.43%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2911'
$c.ClearFormats()
$ws.Range('E8').Value = '  +3.11%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06534'
$c.ClearFormats()
$ws.Range('E9').Value = '  +1.38%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '22.08'
$c.ClearFormats()
$ws.Range('E10').Value = '  +5.49%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07731'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '97.44'
$c.ClearFormats()
$ws.Range('E12').Value = '  +4.48%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.7432'
$c.ClearFormats()
$ws.Range('E13').Value = '  +9.86%  '
$ws.Range('D14').Value = '1.877.80'
$ws.Range('E14').Value = '  +1.05%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.170'
$c.ClearFormats()
$ws.Range('E15').Value = '  +2.72%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '275.50'
$c.ClearFormats()
$ws.Range('E16').Value = '  +3.58%  '
$ws.Range('D17').Value = '30.548.04'
$ws.Range('E17').Value = '  +1.55%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '13.63'
$c.ClearFormats()
$ws.Range('E18').Value = '  +2.64%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.000007588'
$c.ClearFormats()
$ws.Range('E19').Value = '  +0.50%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.9995'
$c.ClearFormats()
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').Value = '2.124.32'
$ws.Range('E21').Value = '  +0.62%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.287'
$c.ClearFormats()
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('E23').Value = '  -0.15%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.221'
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.25%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.335'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.67%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '163.22'
$c.ClearFormats()
$ws.Range('E26').Value = '  -1.23%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '18.95'
$c.ClearFormats()
$ws.Range('E27').Value = '  +2.73%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.950'
$c.ClearFormats()
$ws.Range('E28').Value = '  +3.92%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.371'
$c.ClearFormats()
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('E30').Value = '  +1.75%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.523'
$c.ClearFormats()
$ws.Range('E31').Value = '  +5.08%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.333'
$c.ClearFormats()
$ws.Range('E32').Value = '  +3.57%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.086'
$c.ClearFormats()
$ws.Range('E33').Value = '  +3.14%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.04802'
$c.ClearFormats()
$ws.Range('E34').Value = '  +3.69%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.132'
$c.ClearFormats()
$ws.Range('E35').Value = '  +1.99%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.7044'
$c.ClearFormats()
$ws.Range('E36').Value = '  +2.87%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.719'
$c.ClearFormats()
$ws.Range('E37').Value = '  +0.16%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.01874'
$c.ClearFormats()
$ws.Range('E38').Value = '  +2.61%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.735'
$c.ClearFormats()
$ws.Range('E39').Value = '  +0.63%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '6.349'
$c.ClearFormats()
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('E41').Value = '  +5.79%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '71.49'
$c.ClearFormats()
$ws.Range('E42').Value = '  +1.54%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.4234'
$c.ClearFormats()
$ws.Range('E43').Value = '  +5.13%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.8409'
$c.ClearFormats()
$ws.Range('E44').Value = '  +1.20%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.9994'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.10%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '102.84'
$c.ClearFormats()
$ws.Range('E46').Value = '  +1.13%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '9.310'
$c.ClearFormats()
$ws.Range('E47').Value = '  +2.18%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.116'
$c.ClearFormats()
$ws.Range('E48').Value = '  +2.96%  '
$ws.Range('E49').Value = '  +5.10%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '921.66'
$c.ClearFormats()
$ws.Range('E50').Value = '  +0.78%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.3908'
$c.ClearFormats()
$ws.Range('E51').Value = '  +5.00%  '
